$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point drift on the existing last row (A20)
$ws.Range("A20").Value2 = 45866.87533740741

# Append the new data row (row 21) produced by the scheduled task run
$ws.Range("A20").Copy($ws.Range("A21")) | Out-Null

$ws.Range("A21").Value2 = 45866.91692689954
$ws.Range("B21").Value2 = 2025
$ws.Range("C21").Value2 = 31
$ws.Range("D21").Value2 = 13.11
$ws.Range("E21").Value2 = 89.2
$ws.Range("F21").Value2 = 0
$ws.Range("G21").Value2 = 0
$ws.Range("H21").Value = "-"
$ws.Range("I21").Value2 = 0
$ws.Range("J21").Value = "22:00:22"
